# Issue #52 REST getConfig
# Applies the changes described by the commit to the Issues / Story workbook:
#  - adds three new Issues rows (51, 52, 53) for the new debounce-config /
#    REST get/set config issues
#  - hides rows 32, 45, 46 and 48 (now resolved/not needed in the visible list)
#  - widens the AutoFilter / used range to cover the new rows
#  - nudges the window position and active-cell selection to match the
#    state Excel saved the workbook in afterwards

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# --- new rows -------------------------------------------------------------
$ws.Cells.Item(51, 1).Value = 51
$ws.Cells.Item(51, 2).Value = 3
$ws.Cells.Item(51, 5).Value = "Need to be able to config debounce timeout"
$ws.Cells.Item(51, 7).Value = "General Settings"

$ws.Cells.Item(52, 1).Value = 52
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(52, 4).Value = "arch"
$ws.Cells.Item(52, 5).Value = "REST service for get config"
$ws.Cells.Item(52, 7).Value = "General Settings"

$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 2
$ws.Cells.Item(53, 4).Value = "arch"
$ws.Cells.Item(53, 5).Value = "REST service for set config"
$ws.Cells.Item(53, 7).Value = "General Settings"

$ws.Rows.Item(51).RowHeight = 29

# --- hide resolved/duplicate rows -----------------------------------------
$ws.Rows.Item(32).Hidden = $true
$ws.Rows.Item(45).Hidden = $true
$ws.Rows.Item(46).Hidden = $true
$ws.Rows.Item(48).Hidden = $true

# --- grow the autofilter / used range to include the new rows ------------
$ws.Range("A1:H52").AutoFilter

# --- window / selection bookkeeping ---------------------------------------
$ws.Range("C56").Select()
$excel.ActiveWindow.Left = 7490
